# Adds slide titles for all slides: the deck currently only has the
# title slide (slide 1, "Programming toolkits"). This inserts the
# remaining 11 content slides of the session, each using the
# "Title and Content" layout (ppLayoutText = 2), with only the title
# placeholder populated (content placeholders are left empty, to be
# filled in later), in the deck's natural topic order.

$p = $ppt.ActivePresentation

$titles = @(
    "Object-oriented programming",
    "Scripting languages",
    "Perl",
    "BioPerl",
    "Python",
    "BioPython",
    "Ruby",
    "BioRuby",
    "R",
    "BioConductor",
    "Exercise: a simple script"
)

$index = 2
foreach ($title in $titles) {
    $slide = $p.Slides.Add($index, 2)
    $slide.Shapes.Item(1).TextFrame.TextRange.Text = $title
    $index = $index + 1
}

Write-Host "Slide count: $($p.Slides.Count)"
